$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 7167.1113
$ws.Range("I62").Value = 5501.6665
$ws.Range("K62").Value = 5501.6665
$ws.Range("M62").Value = -4877.6665
$ws.Range("H65").Value = 7167.1113
$ws.Range("I65").Value = 5501.6665
$ws.Range("K65").Value = 27508.3325
$ws.Range("M65").Value = -24388.3325
$ws.Range("H74").Value = 7045.161
$ws.Range("I74").Value = 4169.4287
$ws.Range("K74").Value = 4169.4287
$ws.Range("M74").Value = -3233.4287
$ws.Range("H77").Value = 7045.161
$ws.Range("I77").Value = 4169.4287
$ws.Range("K77").Value = 20847.1435
$ws.Range("M77").Value = -16167.1435
$ws.Range("H106").Value = 3005
$ws.Range("I106").Value = 3005
$ws.Range("K106").Value = 3005
$ws.Range("M106").Value = -2374
$ws.Range("H111").Value = 15873701
$ws.Range("I111").Value = 22223022
$ws.Range("K111").Value = 66669066
$ws.Range("M111").Value = -66665999
$ws.Range("H132").Value = 55561140
$ws.Range("I132").Value = 55561140
$ws.Range("K132").Value = 166683420
$ws.Range("M132").Value = -166680890
$ws.Range("H138").Value = 3360.0278
$ws.Range("I138").Value = 2373.8333
$ws.Range("J138").Value = 3557.2666
$ws.Range("K138").Value = 7121.499899999999
$ws.Range("L138").Value = 10671.7998
$ws.Range("M138").Value = -1981.499899999999
$ws.Range("N138").Value = -20951.7998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H14").Value = 596.5
$ws.Range("I14").Value = 453
$ws.Range("K14").Value = 453
$ws.Range("M14").Value = -278
$ws.Range("H32").Value = 9403.689
$ws.Range("I32").Value = 5454.7554
$ws.Range("J32").Value = 23073.076
$ws.Range("K32").Value = 5454.7554
$ws.Range("L32").Value = 23073.076
$ws.Range("M32").Value = -5167.7554
$ws.Range("N32").Value = -23647.076
$ws.Range("H132").Value = 2662.4736
$ws.Range("I132").Value = 1896
$ws.Range("K132").Value = 5688
$ws.Range("M132").Value = -3158

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3716719
$ws.Range("I94").Value = 6251211.5
$ws.Range("J94").Value = 30184.908
$ws.Range("K94").Value = 6251211.5
$ws.Range("L94").Value = 30184.908
$ws.Range("M94").Value = -6250760.5
$ws.Range("N94").Value = -31086.908
$ws.Range("H139").Value = 111936.5
$ws.Range("J139").Value = 118415.336
$ws.Range("L139").Value = 118415.336
$ws.Range("N139").Value = -128695.336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 2484.4285
$ws.Range("I107").Value = 2314.0833
$ws.Range("K107").Value = 2314.0833
$ws.Range("M107").Value = -394.0832999999998
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()
$ws.Range("H121").Value = 45000
$ws.Range("J121").Value = 45000
$ws.Range("L121").Value = 45000
$ws.Range("N121").Value = -47620
$ws.Range("H132").Value = 34834.555
$ws.Range("I132").Value = 1243.15
$ws.Range("J132").Value = 130810
$ws.Range("K132").Value = 3729.45
$ws.Range("L132").Value = 392430
$ws.Range("M132").Value = -1199.45
$ws.Range("N132").Value = -397490

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H20").Value = 1650
$ws.Range("I20").Value = 800
$ws.Range("J20").Value = 2500
$ws.Range("K20").Value = 2400
$ws.Range("L20").Value = 7500
$ws.Range("M20").Value = -2173
$ws.Range("N20").Value = -7954
$ws.Range("H37").Value = 48700
$ws.Range("J37").Value = 48700
$ws.Range("L37").Value = 146100
$ws.Range("N37").Value = -146324
$ws.Range("H57").Value = 3517.7778
$ws.Range("I57").Value = 1165
$ws.Range("K57").Value = 3495
$ws.Range("M57").Value = -2936
$ws.Range("H98").Value = 2266.5
$ws.Range("J98").Value = 2266.5
$ws.Range("L98").Value = 6799.5
$ws.Range("N98").Value = -9795.5
$ws.Range("H128").Value = 183324.33
$ws.Range("I128").Value = 183324.33
$ws.Range("K128").Value = 549972.99
$ws.Range("M128").Value = -544992.99

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 10031
$ws.Range("I31").Value = 10031
$ws.Range("K31").Value = 10031
$ws.Range("M31").Value = -9739
$ws.Range("H37").Value = 10031
$ws.Range("I37").Value = 10031
$ws.Range("K37").Value = 10031
$ws.Range("M37").Value = -9754
$ws.Range("H80").Value = 22802208
$ws.Range("I80").Value = 32776806
$ws.Range("J80").Value = 3129.5715
$ws.Range("K80").Value = 32776806
$ws.Range("L80").Value = 3129.5715
$ws.Range("M80").Value = -32775808
$ws.Range("N80").Value = -5125.5715
$ws.Range("H83").Value = 22802208
$ws.Range("I83").Value = 32776806
$ws.Range("J83").Value = 3129.5715
$ws.Range("K83").Value = 163884030
$ws.Range("L83").Value = 15647.8575
$ws.Range("M83").Value = -163879038
$ws.Range("N83").Value = -25631.8575
$ws.Range("H132").Value = 3261.2856
$ws.Range("I132").Value = 3266.375
$ws.Range("J132").Value = 3230.75
$ws.Range("K132").Value = 9799.125
$ws.Range("L132").Value = 9692.25
$ws.Range("M132").Value = -7269.125
$ws.Range("N132").Value = -14752.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 15000
$ws.Range("J14").Value = 15000
$ws.Range("L14").Value = 15000
$ws.Range("N14").Value = -15344
$ws.Range("H19").Value = 2418.6
$ws.Range("J19").Value = 2999.5
$ws.Range("L19").Value = 2999.5
$ws.Range("N19").Value = -3339.5
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").ClearContents()
$ws.Range("H141").Value = 117059
$ws.Range("J141").Value = 117059
$ws.Range("L141").Value = 117059
$ws.Range("N141").Value = -127419

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 11117924
$ws.Range("I81").Value = 27781446
$ws.Range("J81").Value = 8909.888999999999
$ws.Range("K81").Value = 55562892
$ws.Range("L81").Value = 17819.778
$ws.Range("M81").Value = -55561831
$ws.Range("N81").Value = -19941.778
$ws.Range("H84").Value = 11117924
$ws.Range("I84").Value = 27781446
$ws.Range("J84").Value = 8909.888999999999
$ws.Range("K84").Value = 277814460
$ws.Range("L84").Value = 89098.88999999998
$ws.Range("M84").Value = -277809156
$ws.Range("N84").Value = -99706.88999999998
$ws.Range("H122").Value = 2030.1714
$ws.Range("I122").Value = 1768.4814
$ws.Range("J122").Value = 2913.375
$ws.Range("K122").Value = 5305.4442
$ws.Range("L122").Value = 8740.125
$ws.Range("M122").Value = -2855.4442
$ws.Range("N122").Value = -13640.125
$ws.Range("H132").Value = 25665762
$ws.Range("I132").Value = 32259970
$ws.Range("K132").Value = 96779910
$ws.Range("M132").Value = -96777380
